$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values per the diff
$ws.Range("H10").Value = 0
$ws.Range("G14").Value = -1
$ws.Range("G15").Value = -1
$ws.Range("G18").Value = -1
$ws.Range("G19").Value = -1
$ws.Range("H24").Value = 0

# Reposition/resize the workbook window to match the saved view state
$win = $wb.Windows.Item(1)
$win.Left = 38280
$win.Top = -120
$win.Width = 21840
$win.Height = 38040

# Update the selection to match the diff (activeCell A10)
$ws.Range("A10").Select() | Out-Null
